$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114-148 down to 115-149.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new data record.
$ws.Range("A114").Value = 6
$ws.Range("B114").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C114").Value = "Metropolitana"
$ws.Range("D114").Value = 44559
$ws.Range("E114").Value = 13
$ws.Range("F114").Value = 100112001
$ws.Range("G114").Value = "Berenjena"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 190
$ws.Range("K114").Value = 7000
$ws.Range("L114").Value = 8000
$ws.Range("M114").Value = 7368
$ws.Range("N114").Value = "$/caja 50 unidades"
$ws.Range("O114").Value = "Región de Arica y Parinacota"
$ws.Range("P114").Value = 147
$ws.Range("Q114").Value = 50
$ws.Range("R114").Value = "Hortaliza"
